$wb = $excel.ActiveWorkbook

# --- 1. Update the absolute path recorded in the workbook (x15ac:absPath) ---
$wb.Path = "D:\Jai Mata Dii\DBS_Automation\ExecutionTestData\6\"

# --- 2. DeviceList sheet: remove the two obsolete device columns (D & E) ---
$ws = $wb.Worksheets.Item("DeviceList")
$ws.Activate()

$ws.Range("D1:E1").EntireColumn.Delete()

# --- 3. Refresh the remaining device columns (B:G) with the current batch data ---
$ws.Range("B1").Value = "APPLE_iPhone11_iOS_14.7.1_28411"
$ws.Range("C1").Value = "APPLE_iPhone7_iOS_13.1.3_316f0"
$ws.Range("D1").Value = "APPLE_iPhone7_iOS_14.1.0_19d9f"
$ws.Range("E1").Value = "APPLE_iPhoneXR_iOS_14.7.1_e2255"
$ws.Range("F1").Value = "APPLE_iPhone8plus_iOS_14.3.0_a0940"
$ws.Range("G1").Value = "APPLE_iPhoneXSMax_iOS_14.0.0_68985"

$ws.Range("B2").Value = "14.7.1"
$ws.Range("C2").Value = "13.1.3"
$ws.Range("D2").Value = "14.1.0"
$ws.Range("E2").Value = "14.7.1"
$ws.Range("F2").Value = "14.3.0"
$ws.Range("G2").Value = "14.0.0"

$ws.Range("B4").Value = "S2021220IUID"
$ws.Range("C4").Value = "S2325476ZUID"
$ws.Range("D4").Value = "S2325485IUID"
$ws.Range("E4").Value = "S2325486GUID"
$ws.Range("F4").Value = "S2325488CUID"
$ws.Range("G4").Value = "S2325490EUID"

$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 5
$ws.Range("G9").Value = 6

$ws.Range("E19").Select()
